$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the "取得日時" (fetched-at) timestamp for all data rows (2-18)
# from 2025-10-08 18:33:29 to 2025-10-09 01:16:19
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-10-09 01:16:19"
}

# Update row 11's price/deadline description: 取引期間 0 日 -> 取引期間 1 日
$ws.Range("D11").Value = "50,000 円 ~ 60,000 円 / 募集期間 3 日、取引期間 1 日"
